# Refresh cryptocurrency price / volume data in the active worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.586.78'
$ws.Range("E2").Value = '  +0.03%  '
$ws.Range("D3").Value = '3.025.79'
$ws.Range("E3").Value = '  +2.30%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '''510.52'
$ws.Range("E5").Value = '  +2.74%  '
$ws.Range("D6").Value = '''140.13'
$ws.Range("E6").Value = '  +3.98%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +1.25%  '
$ws.Range("D9").Value = '''7.13'
$ws.Range("E9").Value = '  -0.68%  '
$ws.Range("E10").Value = '  +1.95%  '
$ws.Range("E11").Value = '  +4.68%  '
$ws.Range("D12").Value = '3.533.50'
$ws.Range("E12").Value = '  +2.07%  '
$ws.Range("E13").Value = '  +0.48%  '
$ws.Range("D14").Value = '''25.29'
$ws.Range("E14").Value = '  -2.51%  '
$ws.Range("E15").Value = '  +3.20%  '
$ws.Range("D16").Value = '56.582.03'
$ws.Range("E16").Value = '  +0.07%  '
$ws.Range("D17").Value = '3.020.63'
$ws.Range("E17").Value = '  +2.36%  '
$ws.Range("E18").Value = '  -1.27%  '
$ws.Range("D19").Value = '''13.09'
$ws.Range("E19").Value = '  +4.76%  '
$ws.Range("E20").Value = '  +3.24%  '
$ws.Range("D21").Value = '''332.24'
$ws.Range("E21").Value = '  +4.64%  '
$ws.Range("E22").Value = '  -0.30%  '
$ws.Range("E23").Value = '  +3.17%  '
$ws.Range("D24").Value = '''64.70'
$ws.Range("E24").Value = '  +3.56%  '
$ws.Range("D25").Value = '3.142.07'
$ws.Range("E25").Value = '  +2.32%  '
$ws.Range("D26").Value = '''0.166'
$ws.Range("E26").Value = '  +2.49%  '
$ws.Range("D27").Value = '''0.999'
$ws.Range("E27").Value = '  +0.17%  '
$ws.Range("D28").Value = '0.0₃0940'
$ws.Range("E28").Value = '  +8.02%  '
$ws.Range("D29").Value = '''6.40'
$ws.Range("E29").Value = '  -1.44%  '
$ws.Range("D30").Value = '''6.79'
$ws.Range("E30").Value = '  -2.93%  '
$ws.Range("E31").Value = '  +2.49%  '
$ws.Range("E32").Value = '  +3.26%  '
$ws.Range("D33").Value = '''20.40'
$ws.Range("E33").Value = '  +2.45%  '
$ws.Range("D34").Value = '''152.73'
$ws.Range("E34").Value = '  +0.06%  '
$ws.Range("D35").Value = '''4.48'
$ws.Range("E35").Value = '  -0.15%  '
$ws.Range("D36").Value = '''27.00'
$ws.Range("E36").Value = '  +13.26%  '
$ws.Range("D37").Value = '''5.82'
$ws.Range("E37").Value = '  +2.33%  '
$ws.Range("D38").Value = '''1.23'
$ws.Range("E38").Value = '  +1.32%  '
$ws.Range("E39").Value = '  +1.17%  '
$ws.Range("D40").Value = '3.061.28'
$ws.Range("E40").Value = '  +2.46%  '
$ws.Range("D41").Value = '''36.51'
$ws.Range("E41").Value = '  -2.17%  '
$ws.Range("D42").Value = '''0.999'
$ws.Range("E42").Value = '  +0.18%  '
$ws.Range("E43").Value = '  +3.18%  '
$ws.Range("D44").Value = '''0.657'
$ws.Range("E44").Value = '  +2.59%  '
$ws.Range("D45").Value = '2.202.59'
$ws.Range("E45").Value = '  +2.42%  '
$ws.Range("E46").Value = '  -0.43%  '
$ws.Range("E47").Value = '  +4.56%  '
$ws.Range("B48").Value = 'ONDO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D48").Value = '''0.927'
$ws.Range("E48").Value = '  +0.73%  '
$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").Value = '''5.85'
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("D50").Value = '''19.70'
$ws.Range("E50").Value = '  +3.51%  '
$ws.Range("E51").Value = '  +0.01%  '
